# Applies the data corrections described in the commit diff to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update individual cell values (ASTHMA=A, OBESITY=B, Died or Recovered=C)
$ws.Range("B2").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("A8").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 2
$ws.Range("A17").Value = 1
$ws.Range("C17").Value = 1

# Move the active cell selection to match the saved view state (I22)
$ws.Activate()
$ws.Range("I22").Select()
